$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.626907116734944)
    3  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    5  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 1, 12.59312877619104)
    6  = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.271902929317955)
    7  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.997429241610044)
    8  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    10 = @(0.04763786555579896, 0.04240448674262143, 0.8054896365839992, 8.660232485948974, 0, 9.555764474831394)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
